$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 5.990141067415806
$ws.Range("D2").Value = 10.01800681850788
$ws.Range("E2").Value = 14.22079268370544
$ws.Range("F2").Value = 29.3093987828187
$ws.Range("G2").Value = 28.34028732308941
$ws.Range("H2").Value = 14.21447172672921
$ws.Range("I2").Value = 19.68792378048687
$ws.Range("J2").Value = 10.09333172774097
$ws.Range("K2").Value = 9.524941405618048
$ws.Range("M2").Value = 14.79781808505996
$ws.Range("N2").Value = 18.37744084083419
$ws.Range("O2").Value = 21.57600693860355

# Row 3
$ws.Range("B3").Value = 5.865021261058802
$ws.Range("D3").Value = 9.998942564881139
$ws.Range("E3").Value = 14.22879470607172
$ws.Range("F3").Value = 29.34729312179984
$ws.Range("G3").Value = 28.36256298529009
$ws.Range("H3").Value = 14.25318810934908
$ws.Range("I3").Value = 19.78030965945597
$ws.Range("J3").Value = 10.11573231375795
$ws.Range("K3").Value = 9.214119897451168
$ws.Range("M3").Value = 14.68570915092243
$ws.Range("N3").Value = 18.42418761653918
$ws.Range("O3").Value = 21.63146523230607

# Row 4
$ws.Range("B4").Value = 5.787732293025661
$ws.Range("D4").Value = 9.988827213803523
$ws.Range("E4").Value = 14.23611692098226
$ws.Range("F4").Value = 29.37711657222443
$ws.Range("G4").Value = 28.3846496307545
$ws.Range("H4").Value = 14.27905110460262
$ws.Range("I4").Value = 19.83996645434897
$ws.Range("J4").Value = 10.13062190614265
$ws.Range("K4").Value = 9.01855894178513
$ws.Range("M4").Value = 14.61847179748038
$ws.Range("N4").Value = 18.45447966996254
$ws.Range("O4").Value = 21.66979989406646

# Row 5
$ws.Range("B5").Value = 5.756164012315336
$ws.Range("D5").Value = 9.985108311902694
$ws.Range("E5").Value = 14.23970732755989
$ws.Range("F5").Value = 29.39091717012001
$ws.Range("G5").Value = 28.39576141426966
$ws.Range("H5").Value = 14.29011635719921
$ws.Range("I5").Value = 19.86501617093034
$ws.Range("J5").Value = 10.13697539985037
$ws.Range("K5").Value = 8.937796469416421
$ws.Range("M5").Value = 14.59149701020788
$ws.Range("N5").Value = 18.46722449644725
$ws.Range("O5").Value = 21.68649715832801

# Row 6
$ws.Range("B6").Value = 5.750919077279537
$ws.Range("D6").Value = 9.984515233963682
$ws.Range("E6").Value = 14.24034016516826
$ws.Range("F6").Value = 29.39330819255357
$ws.Range("G6").Value = 28.39773391322539
$ws.Range("H6").Value = 14.29198549747954
$ws.Range("I6").Value = 19.86922035356769
$ws.Range("J6").Value = 10.13804766755618
$ws.Range("K6").Value = 8.924324891200042
$ws.Range("M6").Value = 14.5870441980967
$ws.Range("N6").Value = 18.46936498560999
$ws.Range("O6").Value = 21.6893346443744

# Row 7
$ws.Range("B7").Value = 5.787306785499251
$ws.Range("D7").Value = 9.988775422725995
$ws.Range("E7").Value = 14.23616288566801
$ws.Range("F7").Value = 29.37729602420458
$ws.Range("G7").Value = 28.38479094490394
$ws.Range("H7").Value = 14.27919820513749
$ws.Range("I7").Value = 19.84030128845661
$ws.Range("J7").Value = 10.13070643369974
$ws.Range("K7").Value = 9.017473921962251
$ws.Range("M7").Value = 14.61810625526644
$ws.Range("N7").Value = 18.45464992803956
$ws.Range("O7").Value = 21.67002072622202

# Row 8
$ws.Range("B8").Value = 5.947122512985373
$ws.Range("D8").Value = 10.01110549908329
$ws.Range("E8").Value = 14.22305233721869
$ws.Range("F8").Value = 29.32110339651351
$ws.Range("G8").Value = 28.34622102453295
$ws.Range("H8").Value = 14.22738720313839
$ws.Range("I8").Value = 19.71917117240021
$ws.Range("J8").Value = 10.10081996295362
$ws.Range("K8").Value = 9.418815270675692
$ws.Range("M8").Value = 14.75884535352204
$ws.Range("N8").Value = 18.39322981150938
$ws.Range("O8").Value = 21.59423919071475

# Row 9
$ws.Range("B9").Value = 6.25495318028967
$ws.Range("D9").Value = 10.06735455980293
$ws.Range("E9").Value = 14.21641494112733
$ws.Range("F9").Value = 29.26296691922427
$ws.Range("G9").Value = 28.33740041856774
$ws.Range("H9").Value = 14.14237332095875
$ws.Range("I9").Value = 19.50480605907694
$ws.Range("J9").Value = 10.05120832615564
$ws.Range("K9").Value = 10.16384257856215
$ws.Range("M9").Value = 15.0464678826535
$ws.Range("N9").Value = 18.2853541445532
$ws.Range("O9").Value = 21.47967550116292

# Row 10
$ws.Range("B10").Value = 6.475333489929419
$ws.Range("D10").Value = 10.11604899226657
$ws.Range("E10").Value = 14.22309999536661
$ws.Range("F10").Value = 29.25200775081799
$ws.Range("G10").Value = 28.371693153319
$ws.Range("H10").Value = 14.09002042072116
$ws.Range("I10").Value = 19.3613111736722
$ws.Range("J10").Value = 10.02022300572103
$ws.Range("K10").Value = 10.68015522186134
$ws.Range("M10").Value = 15.2634422125421
$ws.Range("N10").Value = 18.2137019101826
$ws.Range("O10").Value = 21.41633455622848

# Row 11
$ws.Range("B11").Value = 6.57387029023251
$ws.Range("D11").Value = 10.1397471251554
$ws.Range("E11").Value = 14.22863357807534
$ws.Range("F11").Value = 29.25390831001066
$ws.Range("G11").Value = 28.3961268231645
$ws.Range("H11").Value = 14.0683970458495
$ws.Range("I11").Value = 19.29904456922806
$ws.Range("J11").Value = 10.00730931744389
$ws.Range("K11").Value = 10.90733069837851
$ws.Range("M11").Value = 15.36306563410422
$ws.Range("N11").Value = 18.18274432749805
$ws.Range("O11").Value = 21.39205508547902

# Row 12
$ws.Range("B12").Value = 6.610901351156688
$ws.Range("D12").Value = 10.14893839212616
$ws.Range("E12").Value = 14.23108550153138
$ws.Range("F12").Value = 29.25561625734829
$ws.Range("G12").Value = 28.40664529667762
$ws.Range("H12").Value = 14.06052397649526
$ws.Range("I12").Value = 19.27589666308975
$ws.Range("J12").Value = 10.00258883567133
$ws.Range("K12").Value = 10.9921805657305
$ws.Range("M12").Value = 15.40089626557505
$ws.Range("N12").Value = 18.17125605210775
$ws.Range("O12").Value = 21.38351395978047

# Row 13
$ws.Range("B13").Value = 6.602939170273296
$ws.Range("D13").Value = 10.14694930827263
$ws.Range("E13").Value = 14.2305416148945
$ws.Range("F13").Value = 29.25520450600444
$ws.Range("G13").Value = 28.40432373513761
$ws.Range("H13").Value = 14.06220556300486
$ws.Range("I13").Value = 19.28086283250175
$ws.Range("J13").Value = 10.00359793496792
$ws.Range("K13").Value = 10.97395997948634
$ws.Range("M13").Value = 15.39274450418056
$ws.Range("N13").Value = 18.1737198319007
$ws.Range("O13").Value = 21.38532438987302

# Row 14
$ws.Range("B14").Value = 6.576922721186033
$ws.Range("D14").Value = 10.14049897151048
$ws.Range("E14").Value = 14.2288281669078
$ws.Range("F14").Value = 29.25402903148684
$ws.Range("G14").Value = 28.39696684228822
$ws.Range("H14").Value = 14.06774300603281
$ws.Range("I14").Value = 19.29713154684549
$ws.Range("J14").Value = 10.00691756207311
$ws.Range("K14").Value = 10.91433529814764
$ws.Range("M14").Value = 15.36617602511215
$ws.Range("N14").Value = 18.18179448096662
$ws.Range("O14").Value = 21.39133931061092

# Row 15
$ws.Range("B15").Value = 6.560949052487229
$ws.Range("D15").Value = 10.13657609612722
$ws.Range("E15").Value = 14.22782499274225
$ws.Range("F15").Value = 29.25343764898103
$ws.Range("G15").Value = 28.39262524680897
$ws.Range("H15").Value = 14.07117590205595
$ws.Range("I15").Value = 19.30715269607851
$ws.Range("J15").Value = 10.00897301501325
$ws.Range("K15").Value = 10.87765830837276
$ws.Range("M15").Value = 15.34991498989665
$ws.Range("N15").Value = 18.18677097773014
$ws.Range("O15").Value = 21.39510868479449

# Row 16
$ws.Range("B16").Value = 6.46885622746423
$ws.Range("D16").Value = 10.11453095639685
$ws.Range("E16").Value = 14.22278835727178
$ws.Range("F16").Value = 29.25202195702202
$ws.Range("G16").Value = 28.37027377865537
$ws.Range("H16").Value = 14.09147767211036
$ws.Range("I16").Value = 19.36544086292493
$ws.Range("J16").Value = 10.02109070089417
$ws.Range("K16").Value = 10.66514773369612
$ws.Range("M16").Value = 15.25694780790633
$ws.Range("N16").Value = 18.21575794363511
$ws.Range("O16").Value = 21.41801260687794

# Row 17
$ws.Range("B17").Value = 6.411894538065523
$ws.Range("D17").Value = 10.10139958039544
$ws.Range("E17").Value = 14.22033562366549
$ws.Range("F17").Value = 29.25291596434932
$ws.Range("G17").Value = 28.35882181620859
$ws.Range("H17").Value = 14.10449363195293
$ws.Range("I17").Value = 19.40196845318135
$ws.Range("J17").Value = 10.0288269737961
$ws.Range("K17").Value = 10.53275565367577
$ws.Range("M17").Value = 15.20013207692479
$ws.Range("N17").Value = 18.23395933213373
$ws.Range("O17").Value = 21.43322555278812

# Row 18
$ws.Range("B18").Value = 6.378972231985418
$ws.Range("D18").Value = 10.09399267505713
$ws.Range("E18").Value = 14.21915961405016
$ws.Range("F18").Value = 29.25407847030722
$ws.Range("G18").Value = 28.35306677529696
$ws.Range("H18").Value = 14.11218641133384
$ws.Range("I18").Value = 19.42326156831611
$ws.Range("J18").Value = 10.0333879160999
$ws.Range("K18").Value = 10.45588630657056
$ws.Range("M18").Value = 15.16754145365345
$ws.Range("N18").Value = 18.24458245189326
$ws.Range("O18").Value = 21.44240241564153

# Row 19
$ws.Range("B19").Value = 6.367799056299884
$ws.Range("D19").Value = 10.09151003724298
$ws.Range("E19").Value = 14.21880180953929
$ws.Range("F19").Value = 29.25458346845128
$ws.Range("G19").Value = 28.35126119092803
$ws.Range("H19").Value = 14.11482649388307
$ws.Range("I19").Value = 19.43051978011442
$ws.Range("J19").Value = 10.03495128714084
$ws.Range("K19").Value = 10.42973810278563
$ws.Range("M19").Value = 15.15652281308382
$ws.Range("N19").Value = 18.24820575878444
$ws.Range("O19").Value = 21.44558281289339

# Row 20
$ws.Range("B20").Value = 6.417974978094746
$ws.Range("D20").Value = 10.10278237249545
$ws.Range("E20").Value = 14.22057244025711
$ws.Range("F20").Value = 29.25275371138209
$ws.Range("G20").Value = 28.35995483285807
$ws.Range("H20").Value = 14.10308670501404
$ws.Range("I20").Value = 19.39805070935214
$ws.Range("J20").Value = 10.02799192310712
$ws.Range("K20").Value = 10.54692415971565
$ws.Range("M20").Value = 15.20617126702563
$ws.Range("N20").Value = 18.23200581359225
$ws.Range("O20").Value = 21.43156193159552

# Row 21
$ws.Range("B21").Value = 6.584572330177038
$ws.Range("D21").Value = 10.14238773648519
$ws.Range("E21").Value = 14.22932179005424
$ws.Range("F21").Value = 29.2543474949866
$ws.Range("G21").Value = 28.39909342361373
$ws.Range("H21").Value = 14.0661079695596
$ws.Range("I21").Value = 19.29234134461674
$ws.Range("J21").Value = 10.00593790502104
$ws.Range("K21").Value = 10.93188095872242
$ws.Range("M21").Value = 15.37397718597036
$ws.Range("N21").Value = 18.17941639741694
$ws.Range("O21").Value = 21.38955485446652

# Row 22
$ws.Range("B22").Value = 6.691790603844959
$ws.Range("D22").Value = 10.16953615795497
$ws.Range("E22").Value = 14.23711675272469
$ws.Range("F22").Value = 29.26114799814737
$ws.Range("G22").Value = 28.43204876640503
$ws.Range("H22").Value = 14.04377763652314
$ws.Range("I22").Value = 19.22576619233048
$ws.Range("J22").Value = 9.992513060172158
$ws.Range("K22").Value = 11.17658787731852
$ws.Range("M22").Value = 15.48424992662515
$ws.Range("N22").Value = 18.14641377871467
$ws.Range("O22").Value = 21.36590718016561

# Row 23
$ws.Range("B23").Value = 6.634729613387338
$ws.Range("D23").Value = 10.15493263488788
$ws.Range("E23").Value = 14.2327671102482
$ws.Range("F23").Value = 29.25699232317977
$ws.Range("G23").Value = 28.413786781688
$ws.Range("H23").Value = 14.05552763828793
$ws.Range("I23").Value = 19.26106931672539
$ws.Range("J23").Value = 9.999587774803182
$ws.Range("K23").Value = 11.04663408022883
$ws.Range("M23").Value = 15.42534913236799
$ws.Range("N23").Value = 18.16390301002035
$ws.Range("O23").Value = 21.37817985143039

# Row 24
$ws.Range("B24").Value = 6.41522655121791
$ws.Range("D24").Value = 10.10215676794609
$ws.Range("E24").Value = 14.22046464614499
$ws.Range("F24").Value = 29.25282504569656
$ws.Range("G24").Value = 28.35944001393806
$ws.Range("H24").Value = 14.10372212338276
$ws.Range("I24").Value = 19.39982100880243
$ws.Range("J24").Value = 10.02836909672271
$ws.Range("K24").Value = 10.54052092558462
$ws.Range("M24").Value = 15.20344071818395
$ws.Range("N24").Value = 18.23288850430744
$ws.Range("O24").Value = 21.43231271308801

# Row 25
$ws.Range("B25").Value = 6.172527340463445
$ws.Range("D25").Value = 10.05082615953524
$ws.Range("E25").Value = 14.21617463501049
$ws.Range("F25").Value = 29.27311691759902
$ws.Range("G25").Value = 28.33262380660831
$ws.Range("H25").Value = 14.16359654223613
$ws.Range("I25").Value = 19.56033004529237
$ws.Range("J25").Value = 10.06366856476662
$ws.Range("K25").Value = 9.967389852284237
$ws.Range("M25").Value = 14.96755934207973
$ws.Range("N25").Value = 18.31319786243901
$ws.Range("O25").Value = 21.50701505853499
